# Project5 Beoordeling.xlsx - update assessment checkmarks ("v") on the
# "Blad2" worksheet and move the view/selection, matching the authored
# commit ("token meer en excel bla").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad2")
$ws.Activate()

# --- Row 7/8: mark F and G as achieved ("v") ---------------------------
$ws.Range("F7").Value = "v"
$ws.Range("G7").Value = "v"

$ws.Range("F8").Value = "v"
$ws.Range("G8").Value = "v"

# --- Row 9: mark E,F,G,H as achieved ("v") ------------------------------
$ws.Range("E9").Value = "v"
$ws.Range("F9").Value = "v"
$ws.Range("G9").Value = "v"
$ws.Range("H9").Value = "v"

# --- Row 10 & 12: clear the previously-set "v" in column E -------------
$ws.Range("E10").ClearContents()
$ws.Range("E12").ClearContents()

# --- Row 40: mark H as achieved ("v") -- this flips H38's rollup to "V"
$ws.Range("H40").Value = "v"

# --- Rows 80,81,83,84,85,87,89,90: mark H as achieved ("v") -------------
# (this flips H79's rollup formula result from "O" to "V")
$ws.Range("H80").Value = "v"
$ws.Range("H81").Value = "v"
$ws.Range("H83").Value = "v"
$ws.Range("H84").Value = "v"
$ws.Range("H85").Value = "v"
$ws.Range("H87").Value = "v"
$ws.Range("H89").Value = "v"
$ws.Range("H90").Value = "v"

# --- View state: scroll the frozen pane down and move the selection ----
$ws.Range("D70").Select()
$ws.Range("U82").Select()

$wb.Saved = $false
